$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cap nhat tien do - update progress dates for row 8 (F, G, H)
$ws.Range("F8").Value = 43393
$ws.Range("G8").Value = 43393
$ws.Range("H8").Value = 43393

# Update the active selection to H8
$ws.Range("H8").Select()
